# "fix and add name mappings"
#
# Before:
#   B1 = "29-12-2023"
#   A2 = 86213130 (a stray numeric value)
#   B2 = "Болен" (red fill)
#
# After (3x3 status table):
#            A                    B                         C
#   1                             29-12-2023                01-01-2024
#   2        test                 Не ответил (blue)          Почти выздоровел (orange)
#   3        Панченко Иван        Болен (red)                Здоров (green)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the stray number in A2 -> now a text label ---
$ws.Range("A2").Value = "test"

# --- New name row ---
$ws.Range("A3").Value = "Панченко Иван"

# --- Status cells with their highlight fills ---
# Interior.Color takes a packed BGR integer; values below correspond to the
# RGB hex colors used in the target fills (62C6FF, FFB762, FF6262, 62FF97).
$blue   = 16762466  # RGB 62C6FF
$orange = 6469631   # RGB FFB762
$red    = 6447871   # RGB FF6262
$green  = 9961314   # RGB 62FF97

$ws.Range("B2").Interior.Color = $blue
$ws.Range("B2").Value = "Не ответил"

$ws.Range("C2").Interior.Color = $orange
$ws.Range("C2").Value = "Почти выздоровел"

$ws.Range("B3").Interior.Color = $red
$ws.Range("B3").Value = "Болен"

$ws.Range("C3").Interior.Color = $green
$ws.Range("C3").Value = "Здоров"

# --- New date cell C1. Force it to be stored as text (like B1 already is)
# instead of being auto-parsed into a date serial, then drop the number
# format back to Normal so the cell keeps the default (unstyled) look. ---
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "01-01-2024"
$ws.Range("C1").Style = "Normal"

# --- Column widths sized to fit the new content ---
$ws.Columns.Item(1).ColumnWidth = 15.140625
$ws.Columns.Item(2).ColumnWidth = 11.7109375
$ws.Columns.Item(3).ColumnWidth = 18.5703125
